# Fixed failed script in IAM module
# - Renamed the "truid" placeholder/store key to "userid" and propagated the
#   rename into the dependent request path / validation string.
# - Corrected the password-validation error responses so that "errorcode"
#   actually matches the HTTP "status" value instead of always being 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM")

# Keep/refresh the current selection on the IAM sheet as recorded by the author.
$ws.Activate()

# K29: STORE key for the register-neon-account test changed from "truid" to "userid"
$ws.Range("K29").Value = "userid"

# J30: the "account not activated" validation string now references the
# renamed store key (...truid=(OPQA-2779_userid)...)
$ws.Range("J30").Value = "status=412||error_description=Activate Registered account to continue||truid=(OPQA-2779_userid)||errorcode=412"

# D31: the activation-email API path now uses the renamed store key
$ws.Range("D31").Value = "/account/email/(OPQA-2779_userid)/activate"

# J38-J46: password validation responses - errorcode now mirrors status instead of 0
$ws.Range("J38").Value = "status=422||errorcode=422||reason=New password should not match current password"
$ws.Range("J39").Value = "status=422||errorcode=422||reason=New password should not match previous 4 passwords"
$ws.Range("J41").Value = "status=400||errorcode=400||reason=Update request body is missing required parameters"
$ws.Range("J43").Value = "status=422||errorcode=422"
$ws.Range("J44").Value = "status=422||errorcode=422||reason=Password should be at least 8 characters long||reason=Password should contain at least one alphabet character, either upper or lower case"
$ws.Range("J45").Value = "status=422||errorcode=422||reason=Password should have at least 1 numeric character"
$ws.Range("J46").Value = "status=422||errorcode=422||reason=Password should be at least 8 characters long||reason=Password should have at least 1 numeric character"

# Selection moved from K5 to the L2:L68 range.
$ws.Range("L2:L68").Select()
